$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'women jodphurs',
    'women knee high',
    'women knee high compression white',
    'women knee high stockings',
    'women knee shorts',
    'women knee support',
    'women legging shorts',
    'women legginga',
    'women leggings 2x',
    'women leggings 3 pack',
    'women leggings and tights',
    'women leggings compression',
    'women leggings for gym',
    'women leggings gym',
    'women leggings lululemon',
    'women leggings nike',
    'women leggings pack',
    'women leggings shiny',
    'women leggings with designs',
    'women leggings workout',
    'women leggins capri',
    'women long tights',
    'women lower back support',
    'women lycra',
    'women mid length shorts',
    'women nike compression pants',
    'women nike gear',
    'women nike pro',
    'women opaque brown tights',
    'women pant',
    'women pants tight',
    'women pickleball',
    'women plus size compression leggings',
    'women power stretch leggings',
    'women puma clothing',
    'women puma pants',
    'women race gear',
    'women race runners',
    'women rainbow leggings',
    'women rainbows',
    'women reflective clothing',
    'women riding tights',
    'women runing',
    'women running capris',
    'women running clothes',
    'women running cold',
    'women running compression',
    'women running compression tights',
    'women running gear for cold weather',
    'women running leggings',
    'women running pant',
    'women running pants',
    'women running pants for cold weather',
    'women running shorts',
    'women running shorts long',
    'women running thermal tights',
    'women running tights with pocket',
    'women s bike shorts',
    'women s clothes',
    'women s compression pants',
    'women s running gear',
    'women s running leggings',
    'women s warm tights',
    'women sheer tights',
    'women ski pants fitted',
    'women snowboarding gear',
    'women soccer pants',
    'women sport apparel',
    'women sport leg warmers',
    'women sport leggings',
    'women sport pant',
    'women sport tights',
    'women sports gym yoga running fitness leggings pants',
    'women sports pant',
    'women sports pants',
    'women support knee highs',
    'women tear away pants',
    'women tech pants',
    'women tesla compression pants',
    'women thick tights black',
    'women thin stripes thigh high tights',
    'women tighs',
    'women tight',
    'women tight pants',
    'women tight shorts',
    'women tight top',
    'women tightd',
    'women tights',
    'women tights 3 4',
    'women tights and leggings',
    'women tights gym',
    'women tights leggings',
    'women tights long',
    'women tights pack',
    'women tights short',
    'women tights shorts',
    'women tights sports',
    'women tights suit',
    'women trail pants',
    'women training clothes'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
